# Update the two-digit/one-digit division problems throughout the document.
# Each "old" expression occurs exactly once in the document, so a simple
# Find/Replace (wdReplaceOne) for each pair is sufficient. The pair
# 72÷5= -> 35÷6= is applied before 68÷5= -> 72÷5= so that the freshly
# produced "72÷5=" text (from the second pair) is not re-matched by the
# first pair's search.

$d = $word.ActiveDocument

$replacements = @(
    @("60÷4=", "35÷9="),
    @("26÷2=", "61÷5="),
    @("15÷7=", "15÷2="),
    @("72÷5=", "35÷6="),
    @("68÷5=", "72÷5="),
    @("50÷4=", "15÷3="),
    @("53÷3=", "26÷4="),
    @("79÷5=", "74÷7="),
    @("15÷4=", "27÷9="),
    @("44÷5=", "56÷6="),
    @("32÷2=", "18÷2="),
    @("96÷3=", "96÷5="),
    @("28÷5=", "43÷3="),
    @("88÷2=", "19÷8="),
    @("38÷3=", "87÷9="),
    @("41÷7=", "44÷3="),
    @("78÷7=", "59÷2="),
    @("77÷7=", "15÷6="),
    @("18÷8=", "90÷2="),
    @("32÷4=", "11÷6="),
    @("98÷2=", "60÷3="),
    @("68÷7=", "43÷4="),
    @("72÷9=", "72÷2="),
    @("24÷6=", "73÷4="),
    @("33÷3=", "26÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
